# Generate Report for Handoff
# Adds a new source file entry (730c43ea-7b8b-436c-b1ff-838e0662ebf5.md) as
# a new row to the three report sheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$fileId   = "730c43ea-7b8b-436c-b1ff-838e0662ebf5"
$mdName   = "$fileId.md"
$mdPath   = "e2e\$fileId.md"
$zhXlf    = "$fileId.785c7eff2a9e5c08cbcdf2647011e25b545293f5.zh-cn.xlf"
$deXlf    = "$fileId.785c7eff2a9e5c08cbcdf2647011e25b545293f5.de-de.xlf"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A9").Value = $mdName
$wsOverview.Range("B9").Value = $mdPath
$wsOverview.Range("B9").Style = "HyperLink"
$wsOverview.Range("C9").Value = ".md"
$wsOverview.Range("D9").Value = ""
$wsOverview.Range("E9").Value = "Ready for handoff"
$wsOverview.Range("F9").Value = "Ready for handoff"
$wsOverview.Range("G9").Value = "2016-08-25 20:43:21"
$wsOverview.Range("G9").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B9"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ad387d17d645d96989bf981f81fca11a361c72b1/e2e/730c43ea-7b8b-436c-b1ff-838e0662ebf5.md", "", "", $mdPath) | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A9").Value = $mdName
$wsZhCn.Range("A9").Style = "HyperLink"
$wsZhCn.Range("B9").Value = ".md"
$wsZhCn.Range("C9").Value = "Ready for handoff"
$wsZhCn.Range("D9").Value = "e2e"
$wsZhCn.Range("E9").Value = "ht"
$wsZhCn.Range("F9").Value = "False"
$wsZhCn.Range("G9").Value = $zhXlf
$wsZhCn.Range("H9").Value = "2016-08-25 20:43:16"
$wsZhCn.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I9").Value = ""
$wsZhCn.Range("J9").Value = ""
$wsZhCn.Range("K9").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L9").Value = ""
$wsZhCn.Range("M9").Value = "True"
$wsZhCn.Range("N9").Value = ""
$wsZhCn.Range("O9").Value = "False"
$wsZhCn.Range("P9").Value = ""

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A9"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ad387d17d645d96989bf981f81fca11a361c72b1/e2e/730c43ea-7b8b-436c-b1ff-838e0662ebf5.md", "", "", $mdName) | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A9").Value = $mdName
$wsDeDe.Range("A9").Style = "HyperLink"
$wsDeDe.Range("B9").Value = ".md"
$wsDeDe.Range("C9").Value = "Ready for handoff"
$wsDeDe.Range("D9").Value = "e2e"
$wsDeDe.Range("E9").Value = "ht"
$wsDeDe.Range("F9").Value = "False"
$wsDeDe.Range("G9").Value = $deXlf
$wsDeDe.Range("H9").Value = "2016-08-25 20:43:21"
$wsDeDe.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I9").Value = ""
$wsDeDe.Range("J9").Value = ""
$wsDeDe.Range("K9").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L9").Value = ""
$wsDeDe.Range("M9").Value = "True"
$wsDeDe.Range("N9").Value = ""
$wsDeDe.Range("O9").Value = "False"
$wsDeDe.Range("P9").Value = ""

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A9"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ad387d17d645d96989bf981f81fca11a361c72b1/e2e/730c43ea-7b8b-436c-b1ff-838e0662ebf5.md", "", "", $mdName) | Out-Null
